{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Exact paragraph texts that must be removed from the document.\nconst textsToRemove = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\n\n// Find the index of the \"LOB1004...\" requirement paragraph so we can also\n// remove the single blank paragraph that immediately follows it (the blank\n// paragraph that used to separate the requirements list from the footer).\nlet reqIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"LOB1004: C\u00e1lculo II (Requisito fraco)\") {\n    reqIndex = i;\n    break;\n  }\n}\n\nconst toDelete = [];\nif (reqIndex !== -1 && items[reqIndex + 1] && items[reqIndex + 1].text === \"\") {\n  toDelete.push(items[reqIndex + 1]);\n}\nfor (const p of items) {\n  if (textsToRemove.indexOf(p.text) !== -1) {\n    toDelete.push(p);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"LOB1004: ...\" requirements paragraph so we can also drop the\n# single blank paragraph that used to sit right after it (it separated the\n# requirements list from the page footer that is being removed below).\n$reqIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.Contains(\"LOB1004\")) {\n        $reqIndex = $i\n        break\n    }\n}\n\nif ($reqIndex -ge 1 -and ($reqIndex + 1) -le $d.Paragraphs.Count) {\n    $nextText = $d.Paragraphs.Item($reqIndex + 1).Range.Text\n    $trimmed = $nextText.Trim()\n    if ($trimmed -eq \"\") {\n        $d.Paragraphs.Item($reqIndex + 1).Range.Delete()\n    }\n}\n\n# Walk backwards so deleting a paragraph never shifts the index of the ones\n# still to be inspected.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.Contains(\"Ver no Jupiter\") -or $t.Contains(\"luizeleno@usp.br\")) {\n        $d.Paragraphs.Item($i).Range.Delete()\n    }\n}\n"}
